{"js": "// Highlight quantitative impact metrics (percentages, dollar amounts, large\n// numbers) in specific bullet/paragraph lines by making them bold and\n// colored (#2C3E50). Each target paragraph's run(s) get split so that only\n// the metric substrings carry the new bold+color formatting, leaving the\n// surrounding text formatting untouched (matches the authoring diff).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// For each paragraph whose full text matches one of these entries, bold +\n// color the listed metric substrings (applied in order, left to right).\nconst targets = [\n  {\n    text:\n      \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\",\n    metrics: [\"23%\", \"64%\"],\n  },\n  {\n    text:\n      \"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \u00b14.2% to \u00b12.1%\",\n    metrics: [\"87%\", \"71%\", \"\u00b14.2%\", \"\u00b12.1%\"],\n  },\n  {\n    text: \"\u2022 Wrote RFP and analyzed bids from 1,200 vendors for research platform development\",\n    metrics: [\"1,200\"],\n  },\n  {\n    text:\n      \"\u2022 Created comprehensive meta-analysis framework handling millions of survey responses that became the $400M Polling Consortium Database at The Analyst Institute, now valued at $1B+\",\n    metrics: [\"$400M\", \"$1B\"],\n  },\n  {\n    text: \"\u2022 Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M\",\n    metrics: [\"73.5%\", \"$4.7M\"],\n  },\n  {\n    text: \"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\",\n    metrics: [\"87%\", \"71%\"],\n  },\n];\n\nfor (const para of paragraphs.items) {\n  const target = targets.find((t) => t.text === para.text);\n  if (!target) continue;\n\n  for (const metric of target.metrics) {\n    const found = para.search(metric, { matchCase: true });\n    found.load(\"items\");\n    await context.sync();\n\n    if (found.items.length === 0) continue;\n    const hit = found.items[0];\n    hit.font.bold = true;\n    hit.font.color = \"#2C3E50\";\n    await context.sync();\n  }\n}\n", "ps1": "# Highlight quantitative impact metrics (percentages, dollar amounts, large\n# numbers) in specific bullet/paragraph lines by making them bold and\n# colored (#2C3E50). Each target paragraph is located by its exact original\n# text, then the listed metric substrings inside it are bolded + colored,\n# which splits the paragraph's single run into multiple runs -- matching\n# the authoring diff.\n\n$doc = $word.ActiveDocument\n\n# #2C3E50 as a Word/VBA BGR color value (Font.Color is 0x00BBGGRR).\n$hexR = 0x2C\n$hexG = 0x3E\n$hexB = 0x50\n$highlightColor = $hexB * 65536 + $hexG * 256 + $hexR\n\n$targets = @(\n    @{\n        Text = '\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%'\n        Metrics = @('23%', '64%')\n    },\n    @{\n        Text = '\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \u00b14.2% to \u00b12.1%'\n        Metrics = @('87%', '71%', '\u00b14.2%', '\u00b12.1%')\n    },\n    @{\n        Text = '\u2022 Wrote RFP and analyzed bids from 1,200 vendors for research platform development'\n        Metrics = @('1,200')\n    },\n    @{\n        Text = '\u2022 Created comprehensive meta-analysis framework handling millions of survey responses that became the $400M Polling Consortium Database at The Analyst Institute, now valued at $1B+'\n        Metrics = @('$400M', '$1B')\n    },\n    @{\n        Text = '\u2022 Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M'\n        Metrics = @('73.5%', '$4.7M')\n    },\n    @{\n        Text = '\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%'\n        Metrics = @('87%', '71%')\n    }\n)\n\n$count = $doc.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $doc.Paragraphs.Item($i)\n    $pText = $p.Range.Text.TrimEnd([char]13, [char]7)\n\n    $target = $null\n    foreach ($t in $targets) {\n        if ($t.Text -eq $pText) {\n            $target = $t\n            break\n        }\n    }\n    if ($null -eq $target) {\n        continue\n    }\n\n    $pStart = $p.Range.Start\n    $pEnd = $p.Range.End\n\n    foreach ($metric in $target.Metrics) {\n        $rng = $doc.Range($pStart, $pEnd)\n        $found = $rng.Find.Execute($metric)\n        if ($found) {\n            $rng.Font.Bold = 1\n            $rng.Font.Color = $highlightColor\n        }\n    }\n}\n"}
